$wb = $excel.ActiveWorkbook

# --- Sheet "Version": update Start Time / End Time / Run Time strings ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = "08:18:00 03-16-2017 Thursday Pacific Daylight Time"
$wsVersion.Range("C11").Value = "08:18:36 03-16-2017 Thursday Pacific Daylight Time"
$wsVersion.Range("C12").Value = "0:0:35 0 days"

# --- Sheet "3_3": update BT Idle statistics ---
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = 0.571
$ws33.Range("D3").Value = 0.597
$ws33.Range("E3").Value = 0.545
$ws33.Range("F3").Value = 0.007
$ws33.Range("H3").Value = "0.569652,0.568823,0.57169,0.556466,0.567987,0.57188,0.544584,0.567892,0.567708,0.571601,0.575596,0.574387,0.567702,0.575583,0.568911,0.57978,0.563429,0.575399,0.568911,0.597226,0.568082,0.570589,0.571336,0.578001,0.596948,0.573836,0.573456,0.571424,0.57031,0.568911,0.565569,0.571791,0.568259,0.576153,0.57476,0.570297,0.576812,0.566309,0.572987,0.570576,0.576988,0.564455,0.57243,0.571044,0.569373,0.56334,0.575507,0.579869,0.566493,0.569835,0.566221,0.571778,0.573089,0.595277,0.566765,0.573836,0.568449,0.578843,0.56798,0.572715,0.565942,0.570576,0.566214,0.562043,0.56817,0.56101,0.573374,0.576425,0.562606,0.560467,0.57402,0.573368,0.573646,0.570202,0.575039,0.562036,0.564461,0.567987,0.575691,0.571968,0.576336,0.561859,0.572335,0.569278,0.57281,0.578001,0.566961,0.577729,0.581452,0.563144,0.572722,0.567423,0.574108,0.581825,0.574767,0.576982,0.573836,0.56724,0.574203,0.567335"

$ws33.Range("C4").Value = 9.763999999999999
$ws33.Range("D4").Value = 9.782999999999999
$ws33.Range("E4").Value = 9.743
$ws33.Range("F4").Value = 0.008999999999999999
$ws33.Range("H4").Value = "9.776857,9.767393,9.748175,9.768413,9.774078,9.773134,9.762278,9.754948,9.765247,9.77361,9.771014,9.764506,9.756341,9.771476,9.772394,9.76721,9.751599,9.762278,9.772033,9.764785,9.756334,9.749194,9.761171,9.775743,9.765994,9.754479,9.760981,9.763025,9.76933,9.757075,9.757081,9.765355,9.770356,9.762944,9.752998,9.75579,9.758569,9.770166,9.771565,9.752169,9.754581,9.78289,9.774907,9.76488,9.751232,9.775559,9.779262,9.766001,9.761918,9.749846,9.765335,9.777319,9.77185,9.754853,9.757645,9.765899,9.770628,9.750036,9.757176,9.761531,9.773698,9.773793,9.755689,9.755043,9.768406,9.777414,9.769989,9.742971,9.761544,9.773053,9.764323,9.766083,9.749656,9.766463,9.775471,9.76414,9.748073,9.758094,9.766279,9.772862,9.76043,9.75272,9.756531,9.778528,9.770274,9.761266,9.753922,9.768684,9.77918,9.764038,9.759398,9.759866,9.771565,9.778161,9.760335,9.746972,9.769153,9.777135,9.768413,9.763025"

# --- Sheet "1_8": update Deep Sleep statistics ---
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = 0
$ws18.Range("D3").Value = 0.008
$ws18.Range("E3").Value = -0.007
$ws18.Range("F3").Value = 0.003
$ws18.Range("H3").Value = "0.002247,-0.002723,0.000853,0.00644,0.002271,-0.002199,0.003394,0.006397,0.001421,-0.00133,0.007821,-0.006614,-0.001909,-0.001101,0.000329,0.001433,0.004763,0.000594,-0.000491,-0.002742,0.006717,-0.000799,-0.000509,0.000859,-0.003845,-0.001626,0.002284,0.001686,0.003073,0.001421,-0.002723,-0.000775,-0.00305,-0.003309,-0.00133,-0.001909,-0.001897,-0.004123,-0.001632,-0.005245,0.003918,-0.000805,-0.001595,0.005866,0.004214,-0.003845,2.7e-05,-0.001064,0.000576,0.005034,3.9e-05,-0.003031,-0.001064,0.000872,0.000613,0.00171,0.001975,-0.006337,0.001421,0.001982,0.000884,-0.000238,0.004448,-0.004123,0.002789,0.001667,0.005318,0.004744,0.002833,0.000607,-0.000793,-0.006633,0.006699,0.004491,0.005324,-0.000522,-0.001632,0.00282,0.000323,0.002524,0.000304,-0.004665,0.001982,0.002808,-0.000522,-0.002742,0.000607,-0.004129,-0.006885,0.003357,-0.002723,0.004177,-0.00358,0.00311,-0.002705,-0.001077,-0.002717,-0.000522,-0.001909,0.002278"

$ws18.Range("C4").Value = 0
$ws18.Range("D4").Value = 0.008999999999999999
$ws18.Range("E4").Value = -0.008
$ws18.Range("F4").Value = 0.003
$ws18.Range("H4").Value = "-0.000781,-0.002433,0.000607,0.002253,0.003085,-0.001064,-0.003031,-0.001613,-0.002168,0.003924,0.001155,0.00586,3.9e-05,-0.007736,-0.000781,0.001704,0.003634,0.000588,0.000896,0.000607,-0.00329,0.002253,-0.000509,-0.000503,-0.000805,0.002518,0.000298,0.001408,-0.001089,0.005022,0.001988,-0.001626,-0.005523,-0.003858,-0.004092,-0.000503,-0.00133,0.000341,-0.001058,-0.002181,0.002259,0.001427,0.001433,0.003357,0.006409,-0.005233,-0.001897,-0.00276,0.007796,-0.005504,-0.00133,0.002278,0.005305,0.001433,-0.003309,-0.000787,0.003918,0.002555,0.004479,0.004195,0.001994,-0.001632,0.003375,-0.004678,-0.001626,0.006138,-0.000781,-0.002994,-0.001064,-0.001348,0.003381,0.001125,0.000582,-0.000263,0.002241,0.00253,-0.001909,0.002006,-0.000787,0.005028,0.002808,-0.003272,-0.003864,0.002839,0.007001,-0.002446,3.9e-05,0.003634,0.004226,-0.006059,0.009486,-0.0058,0.001982,0.00089,-0.002181,0.000853,0.007537,-0.002458,0.001982,0.0006"
